$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4: ba/Appreciation -> sd/Statement-non-opinion
$ws.Range("I4").Value = "sd"
$ws.Range("J4").Value = "Statement-non-opinion"

# Row 17: b/Acknowledge (Backchannel) -> %/Uninterpretable
$ws.Range("I17").Value = "%"
$ws.Range("J17").Value = "Uninterpretable"

# Row 22: aa/Agree/Accept -> sd/Statement-non-opinion
$ws.Range("I22").Value = "sd"
$ws.Range("J22").Value = "Statement-non-opinion"

# Row 30: sd/Statement-non-opinion -> b/Acknowledge (Backchannel)
$ws.Range("I30").Value = "b"
$ws.Range("J30").Value = "Acknowledge (Backchannel)"

# Row 31: b/Acknowledge (Backchannel) -> aa/Agree/Accept
$ws.Range("I31").Value = "aa"
$ws.Range("J31").Value = "Agree/Accept"

# Row 34: aa/Agree/Accept -> sd/Statement-non-opinion
$ws.Range("I34").Value = "sd"
$ws.Range("J34").Value = "Statement-non-opinion"

$wb.Save()
